$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.8503067892312051
$ws.Range("B3").Value = 0.9999995761138796
$ws.Range("B4").Value = 22.34493927491718

$ws.Rows.Item(5).Delete()
